# Auto-generated edit script: updates leve-profit computed columns (H:N)
# across multiple sheets to reflect refreshed market-price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 250
$ws.Range("I6").Value = 250
$ws.Range("K6").Value = 750
$ws.Range("M6").Value = -638
$ws.Range("H76").Value = 3008.3333
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 3100
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 3100
$ws.Range("M76").Value = -2685
$ws.Range("N76").Value = -3730
$ws.Range("H79").Value = 3008.3333
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 3100
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 3100
$ws.Range("M79").Value = -1908
$ws.Range("N79").Value = -5284
$ws.Range("H112").Value = 1432.6666
$ws.Range("J112").Value = 1432.6666
$ws.Range("L112").Value = 4297.9998
$ws.Range("N112").Value = -6513.9998
$ws.Range("H131").Value = 125000760
$ws.Range("I131").Value = 142857580
$ws.Range("J131").Value = 3000
$ws.Range("K131").Value = 428572740
$ws.Range("L131").Value = 9000
$ws.Range("M131").Value = -428567700
$ws.Range("N131").Value = -19080
$ws.Range("H132").Value = 31747246
$ws.Range("I132").Value = 3269224.5
$ws.Range("K132").Value = 9807673.5
$ws.Range("M132").Value = -9805143.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H48").Value = 99342
$ws.Range("J48").Value = 99342
$ws.Range("L48").Value = 99342
$ws.Range("N48").Value = -100110
$ws.Range("H52").Value = 18500
$ws.Range("J52").Value = 18500
$ws.Range("L52").Value = 18500
$ws.Range("N52").Value = -19136
$ws.Range("H61").Value = 1292.525
$ws.Range("I61").Value = 1183.5143
$ws.Range("J61").Value = 2055.6
$ws.Range("K61").Value = 1183.5143
$ws.Range("L61").Value = 2055.6
$ws.Range("M61").Value = -971.5143
$ws.Range("N61").Value = -2479.6
$ws.Range("H74").Value = 62411.277
$ws.Range("I74").Value = 78229.69500000001
$ws.Range("J74").Value = 21283.4
$ws.Range("K74").Value = 78229.69500000001
$ws.Range("L74").Value = 21283.4
$ws.Range("M74").Value = -77355.69500000001
$ws.Range("N74").Value = -23031.4
$ws.Range("H77").Value = 62411.277
$ws.Range("I77").Value = 78229.69500000001
$ws.Range("J77").Value = 21283.4
$ws.Range("K77").Value = 391148.475
$ws.Range("L77").Value = 106417
$ws.Range("M77").Value = -386780.475
$ws.Range("N77").Value = -115153
$ws.Range("H102").Value = 1603.2084
$ws.Range("I102").Value = 1258.4667
$ws.Range("J102").Value = 2177.7778
$ws.Range("K102").Value = 1258.4667
$ws.Range("L102").Value = 2177.7778
$ws.Range("M102").Value = 363.5333000000001
$ws.Range("N102").Value = -5421.7778
$ws.Range("H108").Value = 304438.28
$ws.Range("J108").Value = 304438.28
$ws.Range("L108").Value = 304438.28
$ws.Range("N108").Value = -312118.28
$ws.Range("H136").Value = 1292.525
$ws.Range("I136").Value = 1183.5143
$ws.Range("J136").Value = 2055.6
$ws.Range("K136").Value = 3550.5429
$ws.Range("L136").Value = 6166.799999999999
$ws.Range("M136").Value = -1000.5429
$ws.Range("N136").Value = -11266.8
$ws.Range("H137").Value = 68000
$ws.Range("J137").Value = 68000
$ws.Range("L137").Value = 68000
$ws.Range("N137").Value = -78200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 33542.855
$ws.Range("J118").Value = 33542.855
$ws.Range("L118").Value = 33542.855
$ws.Range("N118").Value = -36856.855
$ws.Range("H134").Value = 4904.25
$ws.Range("I134").Value = 3180.5652
$ws.Range("J134").Value = 7953.846
$ws.Range("K134").Value = 9541.695599999999
$ws.Range("L134").Value = 23861.538
$ws.Range("M134").Value = -7006.695599999999
$ws.Range("N134").Value = -28931.538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3489.2307
$ws.Range("I58").Value = 3632.7778
$ws.Range("J58").Value = 1766.6666
$ws.Range("K58").Value = 3632.7778
$ws.Range("L58").Value = 1766.6666
$ws.Range("M58").Value = -3429.7778
$ws.Range("N58").Value = -2172.6666
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H134").Value = 27502006
$ws.Range("I134").Value = 3335379.8
$ws.Range("J134").Value = 100001890
$ws.Range("K134").Value = 10006139.4
$ws.Range("L134").Value = 300005670
$ws.Range("M134").Value = -10003604.4
$ws.Range("N134").Value = -300010740
$ws.Range("H136").Value = 3489.2307
$ws.Range("I136").Value = 3632.7778
$ws.Range("J136").Value = 1766.6666
$ws.Range("K136").Value = 10898.3334
$ws.Range("L136").Value = 5299.9998
$ws.Range("M136").Value = -8348.3334
$ws.Range("N136").Value = -10399.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1127.421
$ws.Range("I4").Value = 240.33333
$ws.Range("J4").Value = 1293.75
$ws.Range("K4").Value = 720.99999
$ws.Range("L4").Value = 3881.25
$ws.Range("M4").Value = -608.99999
$ws.Range("N4").Value = -4105.25
$ws.Range("H5").Value = 18519420
$ws.Range("I5").Value = 47619536
$ws.Range("J5").Value = 1163.5454
$ws.Range("K5").Value = 142858608
$ws.Range("L5").Value = 3490.6362
$ws.Range("M5").Value = -142858496
$ws.Range("N5").Value = -3714.6362
$ws.Range("H7").Value = 221.33333
$ws.Range("I7").Value = 264
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 792
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = -680
$ws.Range("N7").Value = -824
$ws.Range("H92").Value = 2000620
$ws.Range("I92").Value = 800
$ws.Range("J92").Value = 3333833.2
$ws.Range("K92").Value = 2400
$ws.Range("L92").Value = 10001499.6
$ws.Range("M92").Value = -1152
$ws.Range("N92").Value = -10003995.6
$ws.Range("H125").Value = 3414.6667
$ws.Range("I125").Value = 1922
$ws.Range("J125").Value = 6400
$ws.Range("K125").Value = 5766
$ws.Range("L125").Value = 19200
$ws.Range("M125").Value = -846
$ws.Range("N125").Value = -29040
$ws.Range("H131").Value = 911.6799999999999
$ws.Range("J131").Value = 911.6799999999999
$ws.Range("L131").Value = 2735.04
$ws.Range("N131").Value = -12815.04
$ws.Range("H132").Value = 5737073.5
$ws.Range("I132").Value = 2780132
$ws.Range("J132").Value = 15875159
$ws.Range("K132").Value = 25021188
$ws.Range("L132").Value = 142876431
$ws.Range("M132").Value = -25018658
$ws.Range("N132").Value = -142881491
$ws.Range("H133").Value = 8331.875
$ws.Range("J133").Value = 9333.333000000001
$ws.Range("L133").Value = 27999.999
$ws.Range("N133").Value = -38119.999
$ws.Range("H134").Value = 6840.838
$ws.Range("I134").Value = 5275.231
$ws.Range("J134").Value = 7688.875
$ws.Range("K134").Value = 15825.693
$ws.Range("L134").Value = 23066.625
$ws.Range("M134").Value = -10755.693
$ws.Range("N134").Value = -33206.625
$ws.Range("H135").Value = 18519420
$ws.Range("I135").Value = 47619536
$ws.Range("J135").Value = 1163.5454
$ws.Range("K135").Value = 428575824
$ws.Range("L135").Value = 10471.9086
$ws.Range("M135").Value = -428573289
$ws.Range("N135").Value = -15541.9086
$ws.Range("H138").Value = 1177.2667
$ws.Range("I138").Value = 1066.1428
$ws.Range("K138").Value = 3198.4284
$ws.Range("M138").Value = 1941.5716
$ws.Range("H139").Value = 1562.4
$ws.Range("I139").Value = 788.38464
$ws.Range("J139").Value = 2999.8572
$ws.Range("K139").Value = 2365.15392
$ws.Range("L139").Value = 8999.571599999999
$ws.Range("M139").Value = 2774.84608
$ws.Range("N139").Value = -19279.5716
$ws.Range("H140").Value = 5115.5
$ws.Range("I140").Value = 953.3333
$ws.Range("K140").Value = 2859.9999
$ws.Range("M140").Value = 2320.0001
$ws.Range("H141").Value = 2559.8333
$ws.Range("I141").Value = 2559.8333
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7679.499899999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -2499.499899999999
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H132").Value = 2505.9302
$ws.Range("I132").Value = 2425.9656
$ws.Range("J132").Value = 2671.5715
$ws.Range("K132").Value = 7277.8968
$ws.Range("L132").Value = 8014.7145
$ws.Range("M132").Value = -4747.8968
$ws.Range("N132").Value = -13074.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 33500
$ws.Range("J64").Value = 33500
$ws.Range("L64").Value = 33500
$ws.Range("N64").Value = -33950
$ws.Range("H67").Value = 33500
$ws.Range("J67").Value = 33500
$ws.Range("L67").Value = 33500
$ws.Range("N67").Value = -35060
$ws.Range("H132").Value = 2606.2307
$ws.Range("I132").Value = 2417.6492
$ws.Range("J132").Value = 3949.875
$ws.Range("K132").Value = 7252.9476
$ws.Range("L132").Value = 11849.625
$ws.Range("M132").Value = -4722.9476
$ws.Range("N132").Value = -16909.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1266.86
$ws.Range("I136").Value = 855.85
$ws.Range("J136").Value = 2910.9
$ws.Range("K136").Value = 2567.55
$ws.Range("L136").Value = 8732.700000000001
$ws.Range("M136").Value = -17.55000000000018
$ws.Range("N136").Value = -13832.7
